$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.861.35"
$ws.Range("E2").Value = "  -2.86%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.805.19"
$ws.Range("E3").Value = "  -3.21%  "
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.90"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.604"
$ws.Range("E6").Value = "  -1.67%  "
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "38.98"
$ws.Range("E8").Value = "  -8.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.320"
$ws.Range("E9").Value = "  +2.55%  "
$ws.Range("E10").Value = "  -3.17%  "
$ws.Range("E11").Value = "  -1.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.065.87"
$ws.Range("E12").Value = "  -3.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.801.33"
$ws.Range("E13").Value = "  -3.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.658"
$ws.Range("E14").Value = "  -3.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "10.87"
$ws.Range("E15").Value = "  -6.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.54"
$ws.Range("E16").Value = "  -5.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "34.798.60"
$ws.Range("E17").Value = "  -2.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0778"
$ws.Range("E19").Value = "  -3.40%  "
$ws.Range("E20").Value = "  -3.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.73"
$ws.Range("E21").Value = "  -5.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.60"
$ws.Range("E22").Value = "  -4.46%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.22"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.40"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.73"
$ws.Range("E26").Value = "  -3.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.22"
$ws.Range("E27").Value = "  -4.22%  "
$ws.Range("E28").Value = "  -3.67%  "
$ws.Range("E29").Value = "  +5.92%  "
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.91"
$ws.Range("E33").Value = "  -4.04%  "
$ws.Range("E34").Value = "  -8.39%  "
$ws.Range("E35").Value = "  +3.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.682"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "90.24"
$ws.Range("E37").Value = "  -7.97%  "
$ws.Range("E38").Value = "  +1.83%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0191"
$ws.Range("E39").Value = "  -3.04%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.302.71"
$ws.Range("E40").Value = "  -3.53%  "
$ws.Range("E41").Value = "  -0.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.48"
$ws.Range("E42").Value = "  -4.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.954"
$ws.Range("E43").Value = "  -7.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.19"
$ws.Range("E44").Value = "  -13.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.68"
$ws.Range("E45").Value = "  -5.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.12"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("E47").Value = "  -2.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.985.45"
$ws.Range("E48").Value = "  -2.09%  "
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("E50").Value = "  +7.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "98.63"
$ws.Range("E51").Value = "  -6.07%  "
